$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values: raw data re-captured with reduced ("custom") accuracy ---
# (rounded from 3 decimal places down to 2 decimal places)
$row5 = @{
    "B5"  = 15.73
    "C5"  = 11.7
    "D5"  = 1.07
    "E5"  = 34.47
    "F5"  = 27.87
    "G5"  = 12.32
    "H5"  = 48.76
    "I5"  = 19.2
    "J5"  = 8.46
    "K5"  = 12.42
    "L5"  = 13.82
    "M5"  = 14.75
    "N5"  = 3.88
    "O5"  = 12.41
    "P5"  = 17.59
    "Q5"  = 10.58
    "R5"  = 0.74
    "S5"  = 0.69
    "T5"  = 181.52
    "U5"  = 34.71
    "V5"  = 11.45
    "W5"  = 23.2
    "X5"  = 12.17
    "Y5"  = 1.88
    "Z5"  = 23.7
    "AA5" = 10.12
    "AB5" = 9.02
    "AC5" = 10.59
    "AD5" = 14.53
    "AE5" = 0.53
    "AF5" = 44.33
    "AG5" = 6.38
    "AH5" = 14.32
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# --- Remove the last data row (row 6) entirely; only 5 rows of data remain ---
$ws.Rows("6:6").Delete()

# --- Column AA (27th column) gets one unit narrower ---
$ws.Columns.Item(27).ColumnWidth = 6.17
